$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Energy Equipment & Services(32)"
$ws.Range("B2").Value = 0.5423198209279574
$ws.Range("A3").Value = "Containers & Packaging(12)"
$ws.Range("B3").Value = 0.4271537997120528
$ws.Range("A4").Value = "Electric Utilities(28)"
$ws.Range("B4").Value = 0.4204824388844848
$ws.Range("A5").Value = "Marine(15)"
$ws.Range("B5").Value = 0.4187471659623295
$ws.Range("A6").Value = "Road & Rail(22)"
$ws.Range("B6").Value = 0.4163434739926311
$ws.Range("A7").Value = "Multi-Utilities(18)"
$ws.Range("B7").Value = 0.4101314264982653
$ws.Range("A8").Value = "Auto Components(21)"
$ws.Range("B8").Value = 0.4048105841666986
$ws.Range("A9").Value = "Chemicals(51)"
$ws.Range("B9").Value = 0.4001645687494276
$ws.Range("A10").Value = "Machinery(85)"
$ws.Range("B10").Value = 0.3919310531092248
$ws.Range("A11").Value = "Air Freight & Logistics(11)"
$ws.Range("B11").Value = 0.3843055748893925
$ws.Range("A12").Value = "Insurance(75)"
$ws.Range("B12").Value = 0.3663951961680489
$ws.Range("A13").Value = "Metals & Mining(89)"
$ws.Range("B13").Value = 0.354268143137205
$ws.Range("A14").Value = "Consumer Finance(15)"
$ws.Range("B14").Value = 0.3459232203685237
$ws.Range("A15").Value = "Oil, Gas & Consumable Fuels(122)"
$ws.Range("B15").Value = 0.338065111906787
$ws.Range("A16").Value = "Life Sciences Tools & Services(19)"
$ws.Range("B16").Value = 0.3230085395246081
$ws.Range("A17").Value = "Building Products(23)"
$ws.Range("B17").Value = 0.3174692589469055
$ws.Range("A18").Value = "Banks(246)"
$ws.Range("B18").Value = 0.2848110563497159
$ws.Range("A19").Value = "Trading Companies & Distributors(25)"
$ws.Range("B19").Value = 0.2783591086449008
$ws.Range("A20").Value = "Capital Markets(75)"
$ws.Range("B20").Value = 0.2770544365119251
$ws.Range("A21").Value = "Electrical Equipment(28)"
$ws.Range("B21").Value = 0.2747318155666205
$ws.Range("A22").Value = "Diversified Telecommunication Services(20)"
$ws.Range("B22").Value = 0.2696989538733799
$ws.Range("A23").Value = "Media(42)"
$ws.Range("B23").Value = 0.244822527871412
$ws.Range("A24").Value = "Semiconductors & Semiconductor Equipment(68)"
$ws.Range("B24").Value = 0.2392752732412778
$ws.Range("A25").Value = "Aerospace & Defense(37)"
$ws.Range("B25").Value = 0.2389711178610557
$ws.Range("A26").Value = "Construction & Engineering(20)"
$ws.Range("B26").Value = 0.2188648845285672
$ws.Range("A27").Value = "Commercial Services & Supplies(52)"
$ws.Range("B27").Value = 0.2169784159426315
$ws.Range("A28").Value = "Biotechnology(126)"
$ws.Range("B28").Value = 0.1995255792216599
$ws.Range("A29").Value = "Household Durables(39)"
$ws.Range("B29").Value = 0.1903255716778237
$ws.Range("A30").Value = "IT Services(52)"
$ws.Range("B30").Value = 0.1894672843098727
$ws.Range("A31").Value = "Health Care Providers & Services(46)"
$ws.Range("B31").Value = 0.1860657929027088
$ws.Range("A32").Value = "Pharmaceuticals(48)"
$ws.Range("B32").Value = 0.1803163295256444
$ws.Range("A33").Value = "Hotels, Restaurants & Leisure(50)"
$ws.Range("B33").Value = 0.1795038682429143
$ws.Range("A34").Value = "Thrifts & Mortgage Finance(47)"
$ws.Range("B34").Value = 0.1757322938026915
$ws.Range("A35").Value = "Software(66)"
$ws.Range("B35").Value = 0.1725878662950896
$ws.Range("A36").Value = "Health Care Equipment & Supplies(83)"
$ws.Range("B36").Value = 0.1621290761695298
$ws.Range("A37").Value = "Specialty Retail(58)"
$ws.Range("B37").Value = 0.1314788714608264
$ws.Range("A38").Value = "Professional Services(35)"
$ws.Range("B38").Value = 0.1221762510538849
$ws.Range("A39").Value = "Communications Equipment(45)"
$ws.Range("B39").Value = 0.1031403075966969

# Remove now-unused trailing rows 40-41
$ws.Range("A40:B41").Delete() | Out-Null

